# NIT-9001865297.xlsx - "Actualiza base de datos EC y agrega parte 1 de
# nuevos estado de cuenta"
#
# The "Periodo Mora" / "Valor Mora" pairs for the worker's debt periods are
# updated: the period that used to be listed first (2103 / 10902) moves to
# the bottom of the table, and the period that used to be listed last
# (2101 / 36341) moves to the top. The middle row (2102 / 36341) is
# unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: was 2103 / 10902 -> now 2101 / 36341
$ws.Range("E16").Value = "2101"
$ws.Range("F16").Value = 36341

# Row 17: 2102 / 36341 stays the same (no change needed)

# Row 18: was 2101 / 36341 -> now 2103 / 10902
$ws.Range("E18").Value = "2103"
$ws.Range("F18").Value = 10902
